$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix punctuation in a handful of "Razon social"/"Nombre Fantasia" entries
#    where separating commas should have been periods (typo fix in how the
#    names were scraped/exported).
# ---------------------------------------------------------------------------
$nameFixes = @(
    @{Cell="E56";  Value="ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"},
    @{Cell="E222"; Value="ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"},
    @{Cell="E74";  Value="SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"},
    @{Cell="E153"; Value="SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"},
    @{Cell="E219"; Value="SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"},
    @{Cell="E232"; Value="SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"},
    @{Cell="E107"; Value="MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"},
    @{Cell="E212"; Value="MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"},
    @{Cell="E155"; Value="RICCOTTI. MARIANA EDITH"},
    @{Cell="F174"; Value="MERCANZINI. GASTON ARIEL"}
)

foreach ($fix in $nameFixes) {
    $ws.Range($fix.Cell).Value = $fix.Value
}

# ---------------------------------------------------------------------------
# 2) Reformat the "Importe" column (H2:H286). Values were scraped/exported as
#    text using "." thousands separators and "," as the decimal mark, e.g.
#    "20.240,00". The fix drops the thousands separator and swaps the comma
#    for a period, producing plain text like "20240.00" (still text, not a
#    real number -- this column has always held its amounts as strings).
#    Each entry below is "row|newText".
# ---------------------------------------------------------------------------
$importeFixes = @(
    "2|20240.00",
    "3|22440.00",
    "4|6450.00",
    "5|27110.00",
    "6|119480.00",
    "7|15840.00",
    "8|15400.00",
    "9|208900.00",
    "10|177000.00",
    "11|1390.00",
    "12|699.00",
    "13|1325.60",
    "14|1030.00",
    "15|431393.98",
    "16|187200.00",
    "17|522502.20",
    "18|243410.00",
    "19|5132.50",
    "20|302.50",
    "21|1660.80",
    "22|3180.00",
    "23|210.00",
    "24|481339.24",
    "25|52850.00",
    "26|507001.66",
    "27|25493.00",
    "28|7950.00",
    "29|2087.75",
    "30|15809.70",
    "31|24120.00",
    "32|2072.80",
    "33|39981.24",
    "34|1472.00",
    "35|7696.51",
    "36|1080.00",
    "37|5250.00",
    "38|7980.00",
    "39|6890.00",
    "40|190.00",
    "41|71400.00",
    "42|478.47",
    "43|38837.39",
    "44|369.00",
    "45|11710.54",
    "46|1042.64",
    "47|1029.00",
    "48|679.65",
    "49|39.00",
    "50|598.80",
    "51|411.14",
    "52|5910.39",
    "53|12830.98",
    "54|399.00",
    "55|140485.54",
    "56|820.00",
    "57|660.70",
    "58|135.00",
    "59|1143.90",
    "60|12013.00",
    "61|2999.00",
    "62|5158.50",
    "63|3162.19",
    "64|23.00",
    "65|31762.50",
    "66|54200.00",
    "67|49.26",
    "68|3288.55",
    "69|220.00",
    "70|274020.08",
    "71|13210.60",
    "72|1180.00",
    "73|1384.40",
    "74|275.00",
    "75|29707.06",
    "76|6958.87",
    "77|1363.92",
    "78|3645.00",
    "79|1580.28",
    "80|9877.00",
    "81|13095.00",
    "82|12600.00",
    "83|274150.00",
    "84|220.40",
    "85|3228.00",
    "86|34.14",
    "87|308.00",
    "88|799.00",
    "89|89200.71",
    "90|899.50",
    "91|1400.00",
    "92|11400.00",
    "93|24320.00",
    "94|3484.00",
    "95|17110.00",
    "96|48730.00",
    "97|21957.00",
    "98|1980.00",
    "99|127944.00",
    "100|2870.00",
    "101|549.00",
    "102|4350.00",
    "103|1093.00",
    "104|91.09",
    "105|392.00",
    "106|20928.00",
    "107|1345.00",
    "108|80.00",
    "109|87.30",
    "110|9393.00",
    "111|9750.00",
    "112|79.60",
    "113|190.00",
    "114|436.00",
    "115|547661.16",
    "116|371897.61",
    "117|28925.37",
    "118|96744.17",
    "119|1500.00",
    "120|6700.00",
    "121|16.64",
    "122|13570.00",
    "123|231.36",
    "124|14907.58",
    "125|5463.20",
    "126|256.93",
    "127|3389.00",
    "128|300.00",
    "129|10544.00",
    "130|3600.00",
    "131|14962.29",
    "132|59937.00",
    "133|1920.00",
    "134|19773.30",
    "135|1514.25",
    "136|400.00",
    "137|10760.20",
    "138|200.00",
    "139|4810.00",
    "140|6022.00",
    "141|3200.00",
    "142|20471.44",
    "143|88.00",
    "144|7541.50",
    "145|24.90",
    "146|5199.00",
    "147|242500.00",
    "148|98.70",
    "149|9200.00",
    "150|60.00",
    "151|9000.00",
    "152|150.00",
    "153|1020.00",
    "154|780.00",
    "155|4000.00",
    "156|2700.00",
    "157|295.16",
    "158|3350.00",
    "159|8300.00",
    "160|11880.00",
    "161|280.00",
    "162|7870.56",
    "163|3257.00",
    "164|1709.60",
    "165|40667.00",
    "166|25760.00",
    "167|80.00",
    "168|2000.00",
    "169|18000.00",
    "170|1200.00",
    "171|8000.00",
    "172|1440.00",
    "173|2285.00",
    "174|9000.00",
    "175|46750.00",
    "176|2800.00",
    "177|215.86",
    "178|1257.40",
    "179|407.30",
    "180|4287.98",
    "181|300.00",
    "182|4742.00",
    "183|585.00",
    "184|3073.00",
    "185|3448.50",
    "186|13673.68",
    "187|877400.00",
    "188|53698.98",
    "189|4600.00",
    "190|8000.00",
    "191|1500.00",
    "192|1800.00",
    "193|1600.00",
    "194|1657.50",
    "195|1440.00",
    "196|8000.00",
    "197|800.00",
    "198|1000.00",
    "199|6000.00",
    "200|1900.00",
    "201|1200.00",
    "202|5260.00",
    "203|3000.00",
    "204|2880.00",
    "205|6000.00",
    "206|600.00",
    "207|800.00",
    "208|20512.00",
    "209|2600.00",
    "210|1035.00",
    "211|3850.00",
    "212|910.00",
    "213|24800.00",
    "214|4800.00",
    "215|620.00",
    "216|1452.00",
    "217|3840.00",
    "218|24770.00",
    "219|1540.00",
    "220|8500.00",
    "221|563.80",
    "222|9855.00",
    "223|20586.00",
    "224|1258.00",
    "225|2290.00",
    "226|1564.00",
    "227|16.13",
    "228|2950.00",
    "229|6148.49",
    "230|1590.00",
    "231|811.47",
    "232|700.00",
    "233|664.32",
    "234|8220.00",
    "235|8307.07",
    "236|3181.26",
    "237|1000.00",
    "238|150.00",
    "239|1300.00",
    "240|135.00",
    "241|2380.00",
    "242|2346.00",
    "243|345.00",
    "244|401.10",
    "245|1950.00",
    "246|5876.71",
    "247|5133.30",
    "248|33600.00",
    "249|177150.00",
    "250|2864.82",
    "251|85750.00",
    "252|843410.42",
    "253|268.14",
    "254|1560.00",
    "255|360370.00",
    "256|450492.00",
    "257|118000.00",
    "258|370000.00",
    "259|846858.00",
    "260|231276.00",
    "261|32500.00",
    "262|72500.00",
    "263|286737.50",
    "264|502342.80",
    "265|496200.00",
    "266|39000.00",
    "267|44500.00",
    "268|302700.00",
    "269|92780.00",
    "270|200000.00",
    "271|520000.00",
    "272|1160.88",
    "273|4919.70",
    "274|85690.00",
    "275|23000.00",
    "276|900.00",
    "277|24000.00",
    "278|38000.00",
    "279|480.00",
    "280|23600.00",
    "281|44450.00",
    "282|811623.00",
    "283|209900.00",
    "284|218.00",
    "285|933.70",
    "286|27.90"
)

$importeRange = $ws.Range("H2:H286")

# Force text interpretation up front so Excel doesn't coerce these
# numeric-looking strings back into real numbers (which would silently drop
# the trailing zeros / formatting we are trying to write).
$importeRange.NumberFormat = "@"

foreach ($entry in $importeFixes) {
    $parts = $entry.Split("|")
    $row = [int]$parts[0]
    $newValue = $parts[1]
    $ws.Cells.Item($row, 8).Value = $newValue   # column H = 8
}

# Restore the default (no explicit style) on the whole column so only the
# text content changed, matching the original formatting of these cells.
$importeRange.Style = "Normal"

Write-Host "Applied $($nameFixes.Count) name fixes and $($importeFixes.Count) Importe reformats."
